$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.961.16"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.982.11"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "499.21"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.69"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.46"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.358"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "3.502.64"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.83"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "56.994.61"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.09"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "2.990.91"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.81"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.66"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.489"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.67"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("E27").Value = "  -4.93%  "
$ws.Range("D28").Value = "0.0₃0896"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.60"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.13"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.77"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("E32").Value = "  -5.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.22"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.43"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.77"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.11"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0667"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").Value = "3.010.39"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.54"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").Value = "2.200.45"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.39"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("E47").Value = "  -5.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.95"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.16"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.81"
$ws.Range("E51").Value = "  -7.56%  "
